$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column C (between the existing "id"/name col B and the already-styled col D)
# gets a custom width matching the other "18.140625" columns once re-expressed through
# Excel's character-width model (observed empirically: an input ColumnWidth of
# 18.140625 round-trips to a stored sheet width of exactly 19).
$ws.Columns.Item(3).ColumnWidth = 18.140625

# New "Location well" mini-table, mirroring the existing Well_Information table
# (rows 12-14) but for LAS location/county/field header fields.
$ws.Range("D19").Value = "Location well"

$ws.Range("C20").Value = "id_well"
$ws.Range("D20").Value = "COUNTY"
$ws.Range("E20").Value = "LOC"
$ws.Range("F20").Value = "FLD"

$ws.Range("C21").Value = "   15-171-21197"

# Scroll/selection state reflects where the user ended up after adding the table.
$ws.Range("F21").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1

$wb.Save()
